# Generate Report for Handback
# Updates the "latest generated" timestamps for the newly-processed
# c5b0ffec-4d99-4a83-8572-8a84215fda90.md file (row 6 of every sheet).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 6.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-10-24 08:03:39"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) for row 6.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-10-24 08:03:27"
$wsZhCn.Range("K6").Value = "2016-10-24 08:04:07"

# de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) for row 6.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-10-24 08:03:39"
$wsDeDe.Range("K6").Value = "2016-10-24 08:04:25"
